$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.742.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -5.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.808.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "274.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -10.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -7.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.95%  "
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06639"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8316"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07809"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.794.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +30.42%  "
$ws.Range("E15").Value = "  -5.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9991"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007981"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.812.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.950"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.040"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.135"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.652"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "108.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.324"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.192"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08763"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.68%  "
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7245"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -12.35%  "
$ws.Range("E35").Value = "  -7.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.882"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9988"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -8.61%  "
$ws.Range("E39").Value = "  -7.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5165"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -14.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.265"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -15.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9434"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -12.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "112.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.158"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.970"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -14.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  -10.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4546"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.293"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.495"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.72%  "
